# Add rows 72-83 to the "Translation" sheet with new SingleUseId text entries,
# as introduced for the sessionsetup panel config UI + template functions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Each entry: TEXT ID, TYPOGRAPHY NAME, ALIGNMENT, DIRECTION, GB (translated text)
$rows = @(
    @("SingleUseId72", "Default", "Left",   "LTR", "Meas Rate"),
    @("SingleUseId73", "Default", "Center", "LTR", "<value>"),
    @("SingleUseId74", "Default", "Left",   "LTR", "0"),
    @("SingleUseId75", "Default", "Left",   "LTR", "Stamps Number"),
    @("SingleUseId76", "Large",   "Left",   "LTR", "X"),
    @("SingleUseId77", "Default", "Left",   "LTR", "Repeat"),
    @("SingleUseId78", "Default", "Center", "LTR", "<value>"),
    @("SingleUseId79", "Default", "Left",   "LTR", "0"),
    @("SingleUseId80", "Default", "Center", "LTR", "<value> s"),
    @("SingleUseId81", "Default", "Left",   "LTR", "0"),
    @("SingleUseId82", "Default", "Left",   "LTR", "Single"),
    @("SingleUseId83", "Default", "Left",   "LTR", "Continuous")
)

$startRow = 72
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]

    $textCell = $ws.Cells.Item($r, 6)
    $gb = $vals[4]
    if ($gb -eq "0") {
        # Force a purely-numeric-looking translation to be stored as text,
        # matching how this value appears elsewhere in the sheet.
        $textCell.NumberFormat = "@"
    }
    $textCell.Value = $gb
}
